$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh, Tue Sep 10 2024)
# Column D (Price) values are forced to Text format so numeric-looking strings
# (e.g. "1.00", "57.146.57") are preserved exactly as authored, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.146.57'
$ws.Range("E2").Value = '  +3.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.335.61'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.12'
$ws.Range("E5").Value = '  +2.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.92'
$ws.Range("E6").Value = '  +3.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("E8").Value = '  +1.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.349.38'
$ws.Range("E9").Value = '  +1.29%  '
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.27'
$ws.Range("E12").Value = '  +3.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.93'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.765.98'
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.974.80'
$ws.Range("E16").Value = '  +3.54%  '
$ws.Range("E17").Value = '  +2.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.338.98'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.48'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.41'
$ws.Range("E21").Value = '  +3.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.63'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.84'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("E25").Value = '  +8.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.977'
$ws.Range("E26").Value = '  -1.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.94'
$ws.Range("E27").Value = '  +5.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.30'
$ws.Range("E28").Value = '  +13.58%  '
$ws.Range("E29").Value = '  +5.19%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.53'
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.20'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.35'
$ws.Range("E33").Value = '  +1.55%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("E36").Value = '  +1.37%  '
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.04'
$ws.Range("E38").Value = '  +3.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.56'
$ws.Range("E39").Value = '  +7.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.92'
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("E42").Value = '  +4.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '138.64'
$ws.Range("E43").Value = '  +2.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.22'
$ws.Range("E44").Value = '  +5.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '277.16'
$ws.Range("E45").Value = '  +5.68%  '
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0506'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.92'
$ws.Range("E50").Value = '  +7.75%  '
$ws.Range("E51").Value = '  +0.38%  '
